$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data occupies rows 2-49 (4 years x 12 months each, starting 2014-01).
# For each 12-row year block, rotate so the last 3 rows (Oct, Nov, Dec)
# move to the front of the block, followed by the original Jan-Sep rows.
$blockStart = 2
$blockSize = 12
$numBlocks = 4

for ($b = 0; $b -lt $numBlocks; $b++) {
    $start = $blockStart + ($b * $blockSize)

    # Capture the original values for this 12-row block (columns A-D).
    $orig = @()
    for ($i = 0; $i -lt $blockSize; $i++) {
        $r = $start + $i
        $orig += ,@(
            $ws.Cells.Item($r, 1).Value2,
            $ws.Cells.Item($r, 2).Value2,
            $ws.Cells.Item($r, 3).Value2,
            $ws.Cells.Item($r, 4).Value2
        )
    }

    # Build rotated order: last 3 rows first, then the first 9 rows.
    $rotated = @()
    for ($i = $blockSize - 3; $i -lt $blockSize; $i++) { $rotated += ,$orig[$i] }
    for ($i = 0; $i -lt ($blockSize - 3); $i++) { $rotated += ,$orig[$i] }

    for ($i = 0; $i -lt $blockSize; $i++) {
        $r = $start + $i
        $row = $rotated[$i]
        $ws.Cells.Item($r, 1).Value2 = $row[0]
        $ws.Cells.Item($r, 2).Value2 = $row[1]
        $ws.Cells.Item($r, 3).Value2 = $row[2]
        $ws.Cells.Item($r, 4).Value2 = $row[3]
    }
}
